$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Il16", "Kcnj10", "ECs", 3, 1, 4.343956666666666, 13.03187, 0.2551833209483726, 0.2551833209483726, 1, 0.3333333333333333, 0.06440966666666667, 0.193229, 0.02029571630342752, 0.02029571630342752, 0.2797928009144444, 2.51813520823, 0.005179128287334664, 0.005179128287334664),
    @("ECs", "Il16", "Kcnj10", "MuSCs", 3, 1, 4.343956666666666, 13.03187, 0.2551833209483726, 0.2551833209483726, 3, 1, 0.8966493333333334, 2.689948, 0.2825374114598339, 0.2825374114598339, 3.895005849195555, 35.05505264276, 0.07209883494847721, 0.07209883494847721),
    @("ECs", "Il16", "Kcnj10", "Resolving-Mac", 3, 1, 4.343956666666666, 13.03187, 0.2551833209483726, 0.2551833209483726, 3, 1, 2.212500666666667, 6.637502, 0.6971668722367386, 0.6971668722367386, 9.61100702097111, 86.49906318874001, 0.1779053577125608, 0.1779053577125608),
    @("FAPs", "Il16", "Kcnj10", "ECs", 3, 1, 6.440526333333334, 19.321579, 0.3783451411951115, 0.3783451411951115, 1, 0.3333333333333333, 0.06440966666666667, 0.193229, 0.02029571630342752, 0.02029571630342752, 0.4148321542878889, 3.733489388591, 0.00767878565047621, 0.00767878565047621),
    @("FAPs", "Il16", "Kcnj10", "MuSCs", 3, 1, 6.440526333333334, 19.321579, 0.3783451411951115, 0.3783451411951115, 3, 1, 0.8966493333333334, 2.689948, 0.2825374114598339, 0.2825374114598339, 5.774893643099112, 51.974042787892, 0.1068966568316722, 0.1068966568316722),
    @("FAPs", "Il16", "Kcnj10", "Resolving-Mac", 3, 1, 6.440526333333334, 19.321579, 0.3783451411951115, 0.3783451411951115, 3, 1, 2.212500666666667, 6.637502, 0.6971668722367386, 0.6971668722367386, 14.24966880618422, 128.247019255658, 0.2637696987129631, 0.2637696987129631),
    @("MuSCs", "Il16", "Kcnj10", "ECs", 3, 1, 0.3495363333333334, 1.048609, 0.02053331770470026, 0.02053331770470026, 1, 0.3333333333333333, 0.06440966666666667, 0.193229, 0.02029571630342752, 0.02029571630342752, 0.02251351871788889, 0.202621668461, 0.0004167383909027419, 0.0004167383909027419),
    @("MuSCs", "Il16", "Kcnj10", "MuSCs", 3, 1, 0.3495363333333334, 1.048609, 0.02053331770470026, 0.02053331770470026, 3, 1, 0.8966493333333334, 2.689948, 0.2825374114598339, 0.2825374114598339, 0.3134115202591112, 2.820703682332001, 0.005801430432968389, 0.005801430432968389),
    @("MuSCs", "Il16", "Kcnj10", "Resolving-Mac", 3, 1, 0.3495363333333334, 1.048609, 0.02053331770470026, 0.02053331770470026, 3, 1, 2.212500666666667, 6.637502, 0.6971668722367386, 0.6971668722367386, 0.7733493705242224, 6.960144334718001, 0.01431514888082913, 0.01431514888082913),
    @("Resolving-Mac", "Il16", "Kcnj10", "ECs", 3, 1, 5.888867, 17.666601, 0.3459382201518156, 0.3459382201518156, 1, 0.3333333333333333, 0.06440966666666667, 0.193229, 0.02029571630342752, 0.02029571630342752, 0.3792999605143334, 3.413699644629, 0.007021063974713903, 0.007021063974713903),
    @("Resolving-Mac", "Il16", "Kcnj10", "MuSCs", 3, 1, 5.888867, 17.666601, 0.3459382201518156, 0.3459382201518156, 3, 1, 0.8966493333333334, 2.689948, 0.2825374114598339, 0.2825374114598339, 5.280248669638667, 47.52223802674801, 0.09774048924671613, 0.09774048924671613),
    @("Resolving-Mac", "Il16", "Kcnj10", "Resolving-Mac", 3, 1, 5.888867, 17.666601, 0.3459382201518156, 0.3459382201518156, 3, 1, 2.212500666666667, 6.637502, 0.6971668722367386, 0.6971668722367386, 13.02912216341133, 117.262099470702, 0.2411766669303856, 0.2411766669303856)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}